# Apply the commit's changes to the "new_store_best_locations" sheet:
#  - Rename header "Total Weighted Demand" -> "Score"
#  - Round the Score column (A) values to 2 decimal places
#  - Simplify the Coords column (B) tuples from
#       (np.float64(<x>), np.float64(<y>))
#    to
#       (<x rounded to 3 dp>, <y rounded to 3 dp>)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in A1
$ws.Cells.Item(1, 1).Value = "Score"

# Find the last used row in the sheet
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    # --- Column A: round the numeric score to 2 decimal places ---
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($aVal -ne $null) {
        $aRounded = [Math]::Round([double]$aVal, 2)
        $aCell.Value = $aRounded
    }

    # --- Column B: simplify the "(np.float64(x), np.float64(y))" text ---
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($bVal -ne $null -and $bVal -match '\(np\.float64\(([^)]+)\),\s*np\.float64\(([^)]+)\)\)') {
        $x = [double]$matches[1]
        $y = [double]$matches[2]
        $xRounded = [Math]::Round($x, 3)
        $yRounded = [Math]::Round($y, 3)
        $bCell.Value = "(" + $xRounded + ", " + $yRounded + ")"
    }
}
